# "sddp with real rates"
# - rename sheets: goal_data_25 -> goal_data, goal_data -> goal_data_5
# - make the (renamed) "goal_data" sheet the active tab / top-left view
# - add a "cashflows" column (K) on that sheet: K = D - E, with an IRR summary cell
# - drop the stray L27 helper cell that is no longer needed

$wb = $excel.ActiveWorkbook

# --- rename sheets (do the destination-name-free one first to avoid collisions) ---
$sheetGoalData25 = $wb.Worksheets.Item(1)   # currently "goal_data_25"
$sheetGoalData   = $wb.Worksheets.Item(2)   # currently "goal_data"

$sheetGoalData.Name = "goal_data_5"
$sheetGoalData25.Name = "goal_data"

# --- worksheet formerly known as goal_data_25, now "goal_data" ---
$ws = $sheetGoalData25

# cashflows header (copy the header formatting used by the other header cells)
$ws.Range("K1").Value = "cashflows"
$ws.Range("C1").Copy()
$ws.Range("K1").PasteSpecial(-4122)   # xlPasteFormats

# cashflow = inflow (D) - goal cost (E) for each period row
$ws.Range("K2").Formula = "=D2-E2"
$ws.Range("K3:K27").Formula = "=D3-E3"

# summary row: annualised-ish IRR of the cashflow series, shown as a percentage
$ws.Range("K28").Formula = "=IRR(K2:K27)"
$ws.Range("K28").NumberFormat = "0.00%"

# this row's one-off helper value isn't used any more
$ws.Range("L27").ClearContents()

# scroll back to the top and make this the active/selected sheet+cell
$ws.Activate()
$ws.Range("E3").Select()
